$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# --- Edit 1: "ADD – Accumulated Day Degrees by day and cohort" paragraph ---
# Split the single run into three runs with the corrected wording
# ("Day Degrees" -> "Degree Days"), keeping identical run formatting.
$para2 = $tr.Paragraphs(2, 1)
$para2.Characters(1, 18).Text = "ADD – Accumulated "
$para2b = $tr.Paragraphs(2, 1)
$para2b.Characters(19, 12).Text = "Degree Days "
$para2c = $tr.Paragraphs(2, 1)
$para2c.Characters(31, 17).Text = "by day and cohort"

# --- Edit 2: the pi-symbol "Daily Survival Rate" paragraph ---
# Prepend a space onto the third run (leave the pi-symbol run and the
# single-space run that follows it untouched).
$para11 = $tr.Paragraphs(11, 1)
$para11.Characters(3, 21).Text = " – Daily Survival Rate"

# --- Edit 3: merge "B " and "– Daily Eggs Laid" into a single run ---
$para12 = $tr.Paragraphs(12, 1)
$para12.Characters(3, 17).Text = "B – Daily Eggs Laid"
$para12b = $tr.Paragraphs(12, 1)
$para12b.Characters(1, 2).Text = ""

